# Insert a new weekly record at row 191 (pushes existing rows 191..279 down to 192..280)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(191).EntireRow.Insert()

$ws.Range("A191").Value = 5
$ws.Range("B191").Value = "Macroferia Regional de Talca"
$ws.Range("C191").Value = "Maule"
$ws.Range("D191").Value = 44704
$ws.Range("E191").Value = 7
$ws.Range("F191").Value = 100112006
$ws.Range("G191").Value = "Repollo"
$ws.Range("H191").Value = "Crespo record"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 4000
$ws.Range("K191").Value = 1000
$ws.Range("L191").Value = 1000
$ws.Range("M191").Value = 1000
$ws.Range("N191").Value = "$/unidad"
$ws.Range("O191").Value = "Región del Maule"
$ws.Range("P191").Value = 1000
$ws.Range("Q191").Value = 1
$ws.Range("R191").Value = "Hortaliza"

Write-Output "inserted row 191"
